$d = $word.ActiveDocument

# 1) Merge the two runs in paragraph 1 ("Q" + "wqwqwe") into a single run
#    with the new text "VETKA 2222", preserving the existing run formatting.
$d.Content.Find.Execute("Qwqwqwe", $true, $false, $false, $false, $false, $true, 1, $false, "VETKA 2222", 2) | Out-Null

$p1 = $d.Paragraphs.Item(1)

# 2) Move the "_GoBack" bookmark so it sits right after "VETKA 2222" (and
#    before the paragraph mark) in paragraph 1. A bookmark collapsed to a
#    position exactly at the paragraph's text end is ambiguous, so insert a
#    throwaway character first, anchor the bookmark just before it (a clean
#    mid-paragraph position), then remove the throwaway character again.
$endPos = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$endPos.InsertAfter("X")

$p1 = $d.Paragraphs.Item(1)
$bmPos = $d.Range($p1.Range.End - 2, $p1.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null

$p1 = $d.Paragraphs.Item(1)
$tempChar = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$tempChar.Delete()

# 3) Drop the paragraphs that used to hold "2", "Da", "F" and the paragraph
#    that used to carry the bookmark (the bookmark now lives in paragraph 1),
#    leaving the single blank paragraph right after paragraph 1 in place.
$p3 = $d.Paragraphs.Item(3)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($p3.Range.Start, $pLast.Range.End).Delete()
